{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per commit \"done with resume screening\"):\n//  - \"weasdf\"  -> \"PHP Developer\", followed by new \"Industry\" / \"IT-Software\" block\n//  - \"waedzf\"  -> the full Senior Web Developer job-description paragraph\n//  - \"10\"      -> \"20\", followed by new \"salary:\" / \"Experience:\" / \"Required\n//                 Qualification\" heading+body blocks\n//  - \"Express\" -> \"Javascript\", followed by new \"PHP\" / \"Ruby\" list items\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Helper: find the (first) paragraph whose text equals `needle` exactly.\nfunction findParagraph(items, needle) {\n  for (const p of items) {\n    if (p.text === needle) return p;\n  }\n  return null;\n}\n\n// Helper: append a run of new paragraphs after `afterParagraph`. Each entry\n// is [text, builtInStyleOrNull]. Content paragraphs are inserted first (so\n// they inherit the preceding \"Normal\" body paragraph's style with no\n// explicit pStyle), then the heading paragraph is inserted *before* that\n// content paragraph and its style is set explicitly. Returns the last\n// paragraph created (useful as the new cursor for further chaining).\nfunction insertHeadingBodyPairs(afterParagraph, pairs) {\n  let cursor = afterParagraph;\n  for (const [headingText, bodyText] of pairs) {\n    const contentPara = cursor.insertParagraph(bodyText, \"After\");\n    const headingPara = contentPara.insertParagraph(headingText, \"Before\");\n    headingPara.style = \"Heading 1\";\n    cursor = contentPara;\n  }\n  return cursor;\n}\n\n// 1) \"weasdf\" -> \"PHP Developer\" + new Industry / IT-Software block\nconst jobRolePara = findParagraph(paragraphs.items, \"weasdf\");\nif (jobRolePara) {\n  jobRolePara.insertText(\"PHP Developer\", \"Replace\");\n  const industryBody = jobRolePara.insertParagraph(\"IT-Software\", \"After\");\n  const industryHeading = industryBody.insertParagraph(\"Industry\", \"Before\");\n  industryHeading.style = \"Heading 1\";\n}\n\n// 2) \"waedzf\" -> full job description text\nconst jobDescPara = findParagraph(paragraphs.items, \"waedzf\");\nif (jobDescPara) {\n  jobDescPara.insertText(\n    \" We are looking for a Senior Web Developer to build and maintain functional web pages and applications.    Senior Web Developer responsibilities include leading a team of junior developers, refining website specifications and resolving technical issues. To be successful in this role, you should have extensive experience building web pages from scratch and in-depth knowledge of at least one of the following programming languages: Javascript, Ruby or PHP.\",\n    \"Replace\"\n  );\n}\n\n// 3) \"10\" -> \"20\" + new salary / experience / required-qualification blocks\nconst vacanciesPara = findParagraph(paragraphs.items, \"10\");\nif (vacanciesPara) {\n  vacanciesPara.insertText(\"20\", \"Replace\");\n  insertHeadingBodyPairs(vacanciesPara, [\n    [\"salary:\", \"5-10 LPA\"],\n    [\n      \"Experience:\",\n      \"Work experience as a Senior Web Developer  Expertise in at least one programming language, preferably Javascript, Ruby or PHP  Solid knowledge of HTML/CSS  Experience with mockup and UI prototyping tools  Understanding of security practices  Familiarity with network diagnostics tools\"\n    ],\n    [\"Required Qualification\", \"BSc/MSc in Computer Science or relevant field\"]\n  ]);\n}\n\n// 4) \"Express\" -> \"Javascript\" + new PHP / Ruby list items\nconst skillPara = findParagraph(paragraphs.items, \"Express\");\nif (skillPara) {\n  skillPara.insertText(\"Javascript\", \"Replace\");\n  let cursor = skillPara;\n  for (const skill of [\"PHP\", \"Ruby\"]) {\n    const p = cursor.insertParagraph(skill, \"After\");\n    p.style = \"List Number\";\n    cursor = p;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Target change (per commit \"done with resume screening\"):\n#  - \"weasdf\"  -> \"PHP Developer\", followed by new \"Industry\" / \"IT-Software\" block\n#  - \"waedzf\"  -> the full Senior Web Developer job-description paragraph\n#  - \"10\"      -> \"20\", followed by new \"salary:\" / \"Experience:\" / \"Required\n#                 Qualification\" heading+body blocks\n#  - \"Express\" -> \"Javascript\", followed by new \"PHP\" / \"Ruby\" list items\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $exactText) {\n    $idx = 0\n    foreach ($p in $doc.Paragraphs) {\n        $idx = $idx + 1\n        if ($p.Range.Text -eq ($exactText + \"`r\")) {\n            return $idx\n        }\n    }\n    return -1\n}\n\n# Insert a run of heading+body paragraph pairs after the paragraph at\n# $afterIndex (1-based). Each pair is inserted as: body paragraph first\n# (via InsertParagraphAfter -> inherits the preceding body paragraph's\n# \"Normal\" style with no explicit override), then the heading paragraph is\n# inserted *before* that body paragraph and explicitly styled \"Heading 1\".\n# Returns the index of the last paragraph written (the new tail).\nfunction Insert-HeadingBodyPairs($doc, $afterIndex, $pairs) {\n    $cursorIndex = $afterIndex\n    foreach ($pair in $pairs) {\n        $headingText = $pair[0]\n        $bodyText = $pair[1]\n\n        $cursorPara = $doc.Paragraphs.Item($cursorIndex)\n        $cursorPara.Range.InsertParagraphAfter()\n\n        $bodyPara = $doc.Paragraphs.Item($cursorIndex + 1)\n        $bodyPara.Range.Text = $bodyText\n        $bodyPara.Range.InsertParagraphBefore()\n\n        $headingPara = $doc.Paragraphs.Item($cursorIndex + 1)\n        $headingPara.Range.Text = $headingText\n        $headingPara.Style = \"Heading 1\"\n\n        $cursorIndex = $cursorIndex + 2\n    }\n    return $cursorIndex\n}\n\n# 1) \"weasdf\" -> \"PHP Developer\" + new Industry / IT-Software block\n$jobRoleIndex = Find-ParagraphIndex $d \"weasdf\"\nif ($jobRoleIndex -ne -1) {\n    $jobRolePara = $d.Paragraphs.Item($jobRoleIndex)\n    $jobRolePara.Range.Text = \"PHP Developer\"\n    $jobRolePara.Range.InsertParagraphAfter()\n\n    $itPara = $d.Paragraphs.Item($jobRoleIndex + 1)\n    $itPara.Range.Text = \"IT-Software\"\n    $itPara.Range.InsertParagraphBefore()\n\n    $industryPara = $d.Paragraphs.Item($jobRoleIndex + 1)\n    $industryPara.Range.Text = \"Industry\"\n    $industryPara.Style = \"Heading 1\"\n}\n\n# 2) \"waedzf\" -> full job description text\n$jobDescIndex = Find-ParagraphIndex $d \"waedzf\"\nif ($jobDescIndex -ne -1) {\n    $jobDescPara = $d.Paragraphs.Item($jobDescIndex)\n    $jobDescPara.Range.Text = \" We are looking for a Senior Web Developer to build and maintain functional web pages and applications.    Senior Web Developer responsibilities include leading a team of junior developers, refining website specifications and resolving technical issues. To be successful in this role, you should have extensive experience building web pages from scratch and in-depth knowledge of at least one of the following programming languages: Javascript, Ruby or PHP.\"\n}\n\n# 3) \"10\" -> \"20\" + new salary / experience / required-qualification blocks\n$vacanciesIndex = Find-ParagraphIndex $d \"10\"\nif ($vacanciesIndex -ne -1) {\n    $vacanciesPara = $d.Paragraphs.Item($vacanciesIndex)\n    $vacanciesPara.Range.Text = \"20\"\n\n    $pairs = @(\n        , @(\"salary:\", \"5-10 LPA\")\n        , @(\"Experience:\", \"Work experience as a Senior Web Developer  Expertise in at least one programming language, preferably Javascript, Ruby or PHP  Solid knowledge of HTML/CSS  Experience with mockup and UI prototyping tools  Understanding of security practices  Familiarity with network diagnostics tools\")\n        , @(\"Required Qualification\", \"BSc/MSc in Computer Science or relevant field\")\n    )\n    Insert-HeadingBodyPairs $d $vacanciesIndex $pairs | Out-Null\n}\n\n# 4) \"Express\" -> \"Javascript\" + new PHP / Ruby list items\n$skillIndex = Find-ParagraphIndex $d \"Express\"\nif ($skillIndex -ne -1) {\n    $skillPara = $d.Paragraphs.Item($skillIndex)\n    $skillPara.Range.Text = \"Javascript\"\n\n    $cursorIndex = $skillIndex\n    foreach ($skill in @(\"PHP\", \"Ruby\")) {\n        $cursorPara = $d.Paragraphs.Item($cursorIndex)\n        $cursorPara.Range.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($cursorIndex + 1)\n        $newPara.Range.Text = $skill\n        $newPara.Style = \"List Number\"\n        $cursorIndex = $cursorIndex + 1\n    }\n}\n"}
